$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2958.2
$ws.Range("J17").Value = 2958.2
$ws.Range("L17").Value = 8874.599999999999
$ws.Range("N17").Value = -9210.599999999999

$ws.Range("H70").Value = 4244.8335
$ws.Range("I70").Value = 2490
$ws.Range("K70").Value = 7470
$ws.Range("M70").Value = -7200

$ws.Range("H73").Value = 4244.8335
$ws.Range("I73").Value = 2490
$ws.Range("K73").Value = 7470
$ws.Range("M73").Value = -6534

$ws.Range("H87").Value = 34285.5
$ws.Range("J87").Value = 34285.5
$ws.Range("L87").Value = 34285.5
$ws.Range("N87").Value = -36781.5

$ws.Range("H90").Value = 34285.5
$ws.Range("J90").Value = 34285.5
$ws.Range("L90").Value = 102856.5
$ws.Range("N90").Value = -115336.5

$ws.Range("H96").Value = 486.44446
$ws.Range("I96").Value = 481.66666
$ws.Range("J96").Value = 496
$ws.Range("K96").Value = 1444.99998
$ws.Range("L96").Value = 1488
$ws.Range("M96").Value = -71.99998000000005
$ws.Range("N96").Value = -4234

$ws.Range("H100").Value = 6485.8184
$ws.Range("I100").Value = 4067.8
$ws.Range("J100").Value = 8500.833000000001
$ws.Range("K100").Value = 4067.8
$ws.Range("L100").Value = 8500.833000000001
$ws.Range("M100").Value = -3526.8
$ws.Range("N100").Value = -9582.833000000001

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 1257
$ws.Range("I132").Value = 1050.1666
$ws.Range("K132").Value = 3150.4998
$ws.Range("M132").Value = -620.4998000000001

$ws.Range("H141").Value = 4749.857
$ws.Range("I141").Value = 4930.4443
$ws.Range("K141").Value = 14791.3329
$ws.Range("M141").Value = -9611.332900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15251.462
$ws.Range("I32").Value = 11563.862
$ws.Range("J32").Value = 42070.363
$ws.Range("K32").Value = 11563.862
$ws.Range("L32").Value = 42070.363
$ws.Range("M32").Value = -11276.862
$ws.Range("N32").Value = -42644.363

$ws.Range("H61").Value = 6999.778
$ws.Range("I61").Value = 2999.5
$ws.Range("J61").Value = 8142.7144
$ws.Range("K61").Value = 2999.5
$ws.Range("L61").Value = 8142.7144
$ws.Range("M61").Value = -2787.5
$ws.Range("N61").Value = -8566.714400000001

$ws.Range("H63").Value = 7849.5
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 7849.5
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H74").Value = 2627.2058
$ws.Range("I74").Value = 2330.5173
$ws.Range("J74").Value = 4348
$ws.Range("K74").Value = 2330.5173
$ws.Range("L74").Value = 4348
$ws.Range("M74").Value = -1456.5173
$ws.Range("N74").Value = -6096

$ws.Range("H77").Value = 2627.2058
$ws.Range("I77").Value = 2330.5173
$ws.Range("J77").Value = 4348
$ws.Range("K77").Value = 11652.5865
$ws.Range("L77").Value = 21740
$ws.Range("M77").Value = -7284.586499999999
$ws.Range("N77").Value = -30476

$ws.Range("H102").Value = 2811.2727
$ws.Range("I102").Value = 2460.7778
$ws.Range("J102").Value = 4388.5
$ws.Range("K102").Value = 2460.7778
$ws.Range("L102").Value = 4388.5
$ws.Range("M102").Value = -838.7777999999998
$ws.Range("N102").Value = -7632.5

$ws.Range("H122").Value = 2954.75
$ws.Range("I122").Value = 2191.5186
$ws.Range("J122").Value = 5244.4443
$ws.Range("K122").Value = 6574.5558
$ws.Range("L122").Value = 15733.3329
$ws.Range("M122").Value = -4124.5558
$ws.Range("N122").Value = -20633.3329

$ws.Range("H132").Value = 5271.2354
$ws.Range("I132").Value = 4910.387
$ws.Range("K132").Value = 14731.161
$ws.Range("M132").Value = -12201.161

$ws.Range("H136").Value = 6999.778
$ws.Range("I136").Value = 2999.5
$ws.Range("J136").Value = 8142.7144
$ws.Range("K136").Value = 8998.5
$ws.Range("L136").Value = 24428.1432
$ws.Range("M136").Value = -6448.5
$ws.Range("N136").Value = -29528.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H86").Value = 3236.0715
$ws.Range("I86").Value = 3655.125
$ws.Range("K86").Value = 3655.125
$ws.Range("M86").Value = -2532.125

$ws.Range("H89").Value = 3236.0715
$ws.Range("I89").Value = 3655.125
$ws.Range("K89").Value = 18275.625
$ws.Range("M89").Value = -12659.625

$ws.Range("H134").Value = 3492.6904
$ws.Range("I134").Value = 2702.9375
$ws.Range("J134").Value = 6019.9
$ws.Range("K134").Value = 8108.8125
$ws.Range("L134").Value = 18059.7
$ws.Range("M134").Value = -5573.8125
$ws.Range("N134").Value = -23129.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4751.4
$ws.Range("J31").Value = 5904.4116
$ws.Range("L31").Value = 5904.4116
$ws.Range("N31").Value = -6494.4116

$ws.Range("H33").Value = 2000
$ws.Range("I33").Value = 2333.3333
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 2333.3333
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -1954.3333
$ws.Range("N33").Value = -1758

$ws.Range("H34").Value = 4751.4
$ws.Range("J34").Value = 5904.4116
$ws.Range("L34").Value = 5904.4116
$ws.Range("N34").Value = -6308.4116

$ws.Range("H86").Value = 9709.513000000001
$ws.Range("I86").Value = 13601.462
$ws.Range("J86").Value = 1925.6154
$ws.Range("K86").Value = 13601.462
$ws.Range("L86").Value = 1925.6154
$ws.Range("M86").Value = -12478.462
$ws.Range("N86").Value = -4171.6154

$ws.Range("H89").Value = 9709.513000000001
$ws.Range("I89").Value = 13601.462
$ws.Range("J89").Value = 1925.6154
$ws.Range("K89").Value = 68007.31
$ws.Range("L89").Value = 9628.076999999999
$ws.Range("M89").Value = -62391.31
$ws.Range("N89").Value = -20860.077

$ws.Range("H132").Value = 4714.5
$ws.Range("J132").Value = 6233.4
$ws.Range("L132").Value = 18700.2
$ws.Range("N132").Value = -23760.2

$ws.Range("H134").Value = 3432.3547
$ws.Range("I134").Value = 1927.6818
$ws.Range("K134").Value = 5783.0454
$ws.Range("M134").Value = -3248.0454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 3203.3333
$ws.Range("I17").Value = 3644
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 10932
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -10763
$ws.Range("N17").Value = -3338

$ws.Range("H117").Value = 3791.3076
$ws.Range("I117").Value = 2694
$ws.Range("J117").Value = 3990.818
$ws.Range("K117").Value = 8082
$ws.Range("L117").Value = 11972.454
$ws.Range("M117").Value = -4640
$ws.Range("N117").Value = -18856.454

$ws.Range("H131").Value = 4237.08
$ws.Range("I131").Value = 3386.7693
$ws.Range("J131").Value = 5158.25
$ws.Range("K131").Value = 10160.3079
$ws.Range("L131").Value = 15474.75
$ws.Range("M131").Value = -5120.3079
$ws.Range("N131").Value = -25554.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H102").Value = 2688.7666
$ws.Range("I102").Value = 2189.16
$ws.Range("J102").Value = 5186.8
$ws.Range("K102").Value = 2189.16
$ws.Range("L102").Value = 5186.8
$ws.Range("M102").Value = -567.1599999999999
$ws.Range("N102").Value = -8430.799999999999

$ws.Range("H122").Value = 15128.765
$ws.Range("I122").Value = 18399.924
$ws.Range("J122").Value = 4497.5
$ws.Range("K122").Value = 55199.772
$ws.Range("L122").Value = 13492.5
$ws.Range("M122").Value = -52749.772
$ws.Range("N122").Value = -18392.5

$ws.Range("H132").Value = 5134.8335
$ws.Range("I132").Value = 5142.794
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 15428.382
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -12898.382
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3817.5
$ws.Range("I12").Value = 7500
$ws.Range("J12").Value = 135
$ws.Range("K12").Value = 7500
$ws.Range("L12").Value = 135
$ws.Range("M12").Value = -7330
$ws.Range("N12").Value = -475

$ws.Range("H22").Value = 1700.6154
$ws.Range("I22").Value = 887.625
$ws.Range("J22").Value = 3001.4
$ws.Range("K22").Value = 887.625
$ws.Range("L22").Value = 3001.4
$ws.Range("M22").Value = -592.625
$ws.Range("N22").Value = -3591.4

$ws.Range("H27").Value = 1700.6154
$ws.Range("I27").Value = 887.625
$ws.Range("J27").Value = 3001.4
$ws.Range("K27").Value = 887.625
$ws.Range("L27").Value = 3001.4
$ws.Range("M27").Value = -780.625
$ws.Range("N27").Value = -3215.4

$ws.Range("H50").Value = 42487.5
$ws.Range("J50").Value = 42487.5
$ws.Range("L50").Value = 42487.5
$ws.Range("N50").Value = -43761.5

$ws.Range("H132").Value = 5360.909
$ws.Range("I132").Value = 4616.2
$ws.Range("J132").Value = 7688.125
$ws.Range("K132").Value = 13848.6
$ws.Range("L132").Value = 23064.375
$ws.Range("M132").Value = -11318.6
$ws.Range("N132").Value = -28124.375

$ws.Range("H136").Value = 4955.923
$ws.Range("I136").Value = 3505.4736
$ws.Range("K136").Value = 10516.4208
$ws.Range("M136").Value = -7966.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3336.125
$ws.Range("I132").Value = 2259.8462
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 6779.5386
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -4249.5386
$ws.Range("N132").Value = -29060

$ws.Range("H136").Value = 3666.9524
$ws.Range("I136").Value = 3082.5862
$ws.Range("J136").Value = 4970.5386
$ws.Range("K136").Value = 9247.758600000001
$ws.Range("L136").Value = 14911.6158
$ws.Range("M136").Value = -6697.758600000001
$ws.Range("N136").Value = -20011.6158
